$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add role / department / group columns ---
# Copy existing header formatting (bold, centered, wrapped) onto the new header cells
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

$ws.Range("E1").Value = "role"
$ws.Range("F1").Value = "department"
$ws.Range("G1").Value = "group"

# --- Data row (row 2): replace sample user with student_1 record ---
# Copy existing body formatting onto the new F2 cell first
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("A2").Value = "student_1"
$ws.Range("B2").Value = "student_1"
$ws.Range("C2").Value = "stu1@gmail.com"
$ws.Range("D2").Value = "password"
$ws.Range("E2").Value = "student"
$ws.Range("F2").Value = "Computer Science"

# Turn the e-mail address into a mailto hyperlink (Excel auto-applies the
# built-in "Hyperlink" style/font to the cell)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:stu1@gmail.com")

# --- Data row (row 3): blank it out, keep formatting ---
$ws.Range("A3:E3").ClearContents()

# --- Selection cursor ends on H2, matching the source workbook ---
$ws.Range("H2").Select()
